$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the report title in the merged A1:J1 banner.
$ws.Range("A1").Value = "SUSPENSE PAYMENT / ETC. CONTROL SHEET"

# 2. Add a blank labeled cell next to "Account" (row 4) that picks up the
#    existing wrap-text style but left-aligned instead of centered.
$ws.Range("B4").HorizontalAlignment = -4131

# 3. Add the form number to the right of the "Currency" row, using a
#    Times New Roman font (left aligned, inheriting the column's format).
$ws.Range("J5").Value = "Form: III-15"
$ws.Range("J5").Font.Name = "Times New Roman"
